# Auto-generated edit script: apply scheduled-runner price/profit updates
# to the Zodiark_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7592.4
$ws.Range("I69").Value = 6796.9
$ws.Range("J69").Value = 9183.4
$ws.Range("K69").Value = 20390.7
$ws.Range("L69").Value = 27550.2
$ws.Range("M69").Value = -19516.7
$ws.Range("N69").Value = -29298.2
$ws.Range("H72").Value = 7592.4
$ws.Range("I72").Value = 6796.9
$ws.Range("J72").Value = 9183.4
$ws.Range("K72").Value = 61172.1
$ws.Range("L72").Value = 82650.59999999999
$ws.Range("M72").Value = -56804.1
$ws.Range("N72").Value = -91386.59999999999
$ws.Range("H132").Value = 2568.2222
$ws.Range("I132").Value = 2615.2144
$ws.Range("J132").Value = 1910.3334
$ws.Range("K132").Value = 7845.6432
$ws.Range("L132").Value = 5731.0002
$ws.Range("M132").Value = -5315.6432
$ws.Range("N132").Value = -10791.0002
$ws.Range("H137").Value = 5634.72
$ws.Range("I137").Value = 7841.5
$ws.Range("J137").Value = 1711.5555
$ws.Range("K137").Value = 23524.5
$ws.Range("L137").Value = 5134.666499999999
$ws.Range("M137").Value = -20974.5
$ws.Range("N137").Value = -10234.6665
$ws.Range("H138").Value = 1716.63
$ws.Range("J138").Value = 2226.8254
$ws.Range("L138").Value = 6680.476200000001
$ws.Range("N138").Value = -16960.4762
$ws.Range("H141").Value = 5680.4043
$ws.Range("I141").Value = 2859.342
$ws.Range("K141").Value = 8578.026
$ws.Range("M141").Value = -3398.026

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4349.793
$ws.Range("I122").Value = 4349.793
$ws.Range("K122").Value = 13049.379
$ws.Range("M122").Value = -10599.379

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 5909.3335
$ws.Range("I97").Value = 5909.3335
$ws.Range("K97").Value = 5909.3335
$ws.Range("M97").Value = -4918.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1866.3334
$ws.Range("I35").Value = 1799.5
$ws.Range("K35").Value = 1799.5
$ws.Range("M35").Value = -1505.5
$ws.Range("H64").Value = 64799
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 64799
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H69").Value = 61795
$ws.Range("I69").Value = 61795
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 61795
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -61046
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 61795
$ws.Range("I72").Value = 61795
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 185385
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -181641
$ws.Range("N72").ClearContents()
$ws.Range("H74").Value = 69532.2
$ws.Range("J74").Value = 69532.2
$ws.Range("L74").Value = 69532.2
$ws.Range("N74").Value = -71280.2
$ws.Range("H77").Value = 69532.2
$ws.Range("J77").Value = 69532.2
$ws.Range("L77").Value = 208596.6
$ws.Range("N77").Value = -217332.6
$ws.Range("H88").Value = 17085.75
$ws.Range("J88").Value = 17085.75
$ws.Range("L88").Value = 17085.75
$ws.Range("N88").Value = -17897.75
$ws.Range("H91").Value = 17085.75
$ws.Range("J91").Value = 17085.75
$ws.Range("L91").Value = 17085.75
$ws.Range("N91").Value = -19893.75
$ws.Range("H122").Value = 1309.35
$ws.Range("I122").Value = 1353.4615
$ws.Range("K122").Value = 4060.3845
$ws.Range("M122").Value = -1610.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 12507
$ws.Range("J88").Value = 12507
$ws.Range("L88").Value = 37521
$ws.Range("N88").Value = -38377
$ws.Range("H91").Value = 12507
$ws.Range("J91").Value = 12507
$ws.Range("L91").Value = 37521
$ws.Range("N91").Value = -40485
$ws.Range("H93").Value = 2400
$ws.Range("I93").Value = 2400
$ws.Range("K93").Value = 7200
$ws.Range("M93").Value = -5328
$ws.Range("H107").Value = 604.3333
$ws.Range("J107").Value = 266.33334
$ws.Range("L107").Value = 799.0000200000001
$ws.Range("N107").Value = -4639.00002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 185.66667
$ws.Range("I9").Value = 57
$ws.Range("K9").Value = 57
$ws.Range("M9").Value = 113
$ws.Range("H99").Value = 9593.200000000001
$ws.Range("I99").Value = 9593.200000000001
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 9593.200000000001
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -7347.200000000001
$ws.Range("N99").ClearContents()
$ws.Range("H118").Value = 15921.8
$ws.Range("J118").Value = 15921.8
$ws.Range("L118").Value = 15921.8
$ws.Range("N118").Value = -19235.8
$ws.Range("H132").Value = 19862.842
$ws.Range("I132").Value = 19075.166
$ws.Range("K132").Value = 57225.49800000001
$ws.Range("M132").Value = -54695.49800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3999
$ws.Range("I93").Value = 3433.4
$ws.Range("J93").Value = 6827
$ws.Range("K93").Value = 3433.4
$ws.Range("L93").Value = 6827
$ws.Range("M93").Value = -2185.4
$ws.Range("N93").Value = -9323
$ws.Range("H122").Value = 5649.4595
$ws.Range("I122").Value = 3440.4783
$ws.Range("J122").Value = 9278.5
$ws.Range("K122").Value = 10321.4349
$ws.Range("L122").Value = 27835.5
$ws.Range("M122").Value = -7871.4349
$ws.Range("N122").Value = -32735.5
$ws.Range("H132").Value = 4326.923
$ws.Range("I132").Value = 5055.7
$ws.Range("J132").Value = 1897.6666
$ws.Range("K132").Value = 15167.1
$ws.Range("L132").Value = 5692.9998
$ws.Range("M132").Value = -12637.1
$ws.Range("N132").Value = -10752.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1999
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H7").Value = 9187.125
$ws.Range("I7").Value = 9187.125
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 9187.125
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -9074.125
$ws.Range("N7").ClearContents()
$ws.Range("H107").Value = 809.6875
$ws.Range("I107").Value = 796.6923
$ws.Range("J107").Value = 866
$ws.Range("K107").Value = 2390.0769
$ws.Range("L107").Value = 2598
$ws.Range("M107").Value = -470.0769
$ws.Range("N107").Value = -6438
$ws.Range("H122").Value = 15155448
$ws.Range("J122").Value = 3727.7778
$ws.Range("L122").Value = 11183.3334
$ws.Range("N122").Value = -16083.3334
$ws.Range("H136").Value = 6416.075
$ws.Range("I136").Value = 8542.5
$ws.Range("K136").Value = 25627.5
$ws.Range("M136").Value = -23077.5

